$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New invoice line for board components, row 6, columns M:O
$ws.Range("O6").Value = "Oct.31"
$ws.Range("M6").Value = "Board Components (See Component Invoice.xlsx)"
$ws.Range("N6").Value = 64.93

# Wrap text for the new items-ordered cell (matches existing style for that column)
$ws.Range("M6").WrapText = $true
$ws.Range("N6").WrapText = $true
$ws.Range("N6").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

# Row height grows to fit the wrapped text
$ws.Rows(6).RowHeight = 44.25

# Leave the new active selection on M6 (matches workbook's last saved selection)
$ws.Range("M6").Select()
